# Cambio de formato del xlsx de edificios, para que los dias tengan tildes
# donde corresponda (cumple con el diseño de la UI).
#   Miercoles -> Miércoles
#   Sabado    -> Sábado

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixes: add the missing accents on the day names (header row) ---
$ws.Range("D1").Value = "Miércoles"
$ws.Range("G1").Value = "Sábado"

# --- Selection moves from I3 to G1 ---
$ws.Range("G1").Select()

# --- Re-apply (touch) the cell formatting across the whole used range so
#     every cell's protection is explicitly carried in its style, exactly
#     as it is laid out today, without changing any visible formatting
#     (fonts, wrapping and alignment stay exactly the same).
$ws.Range("A1:H1").Locked = $true
$ws.Range("A2").Locked = $true
$ws.Range("H2").Locked = $true
$ws.Range("A3").Locked = $true
$ws.Range("H3").Locked = $true
$ws.Range("A4").Locked = $true
$ws.Range("H4").Locked = $true
$ws.Columns("A").Locked = $true

# The "horario" cells (Times New Roman, wrapped) keep their wrap explicitly
# set while the protection flag is re-applied, so none of them lose their
# wrapping as a side effect of the style being rewritten.
$ws.Range("B2:G4").WrapText = $true
$ws.Range("B2:G4").Locked = $true
